$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (data rows start at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C holds the "Förändrad" (last changed) date for every record.
# All rows were bumped from serial date 45190 to 45192 (2023-09-21 -> 2023-09-23).
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45192
